$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Requirement text that used to live in row 24 is replaced with a new
# (different) requirement.
$newReq1 = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$ws.Range("B24").Value = $newReq1
$ws.Range("C24").Value = $newReq1

# Add a new row 25, mirroring row 24's formatting (styles + 30pt row
# height) to hold an additional requirement line.
$ws.Range("B24:C24").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newReq2 = "LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)`n"
$ws.Range("B25").Value = $newReq2
$ws.Range("C25").Value = $newReq2

$ws.Rows.Item(25).RowHeight = 30
